$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
  "47+29=",
  "19+76=",
  "78+16=",
  "19+57=",
  "75-67=",
  "94-68=",
  "50-7=",
  "53-29=",
  "50-35=",
  "56+8=",
  "56+7=",
  "87-29=",
  "92-4=",
  "29+23=",
  "13+59=",
  "96-67=",
  "41-32=",
  "72-34=",
  "19+53=",
  "30-21=",
  "26+46=",
  "47+4=",
  "59+25=",
  "39+28=",
  "55-27=",
  "84-9=",
  "60-53=",
  "92-76=",
  "51-4=",
  "50-19=",
  "16+35=",
  "9+84=",
  "4+47=",
  "65+8=",
  "40-37=",
  "7+76=",
  "91-54=",
  "49+16=",
  "56-38=",
  "73-8=",
  "37+49=",
  "83-75=",
  "16+27=",
  "54-26=",
  "33-6=",
  "53-27=",
  "55-48=",
  "41-22=",
  "32-8=",
  "50-36=",
  "35+29=",
  "35-18=",
  "76-8=",
  "43-5=",
  "35+19=",
  "40-32=",
  "61-53=",
  "48+13=",
  "39+19=",
  "68+24=",
  "41-28=",
  "54+18=",
  "61-54=",
  "37-28=",
  "62-4=",
  "89+6=",
  "90-11=",
  "29+16=",
  "6+86=",
  "48+26=",
  "72-25=",
  "64-17=",
  "19+27=",
  "60-53=",
  "58+16=",
  "24+7=",
  "27-9=",
  "49+25=",
  "97-89=",
  "51-42=",
  "20-1=",
  "57+8=",
  "27+37=",
  "93-87=",
  "20-5=",
  "67+6=",
  "12-4=",
  "16+7=",
  "62-6=",
  "58+29=",
  "74-47=",
  "47+24=",
  "41-14=",
  "15+28=",
  "72-59=",
  "87+5=",
  "25-9=",
  "27+26=",
  "45+37=",
  "38+9="
)

$cols = 5
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = [int][math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newValues[$i]
}

$d.Save()
